# OCEPROJECT-4609: Testing Page Options and Page Options buttons.
$wb = $excel.ActiveWorkbook

# Create the new "All Pages" sheet by copying "Not Visible Page Options" (sheet 2),
# which already carries the exact column widths / formatting we need and has no
# pageSetup element, placing the copy before the first sheet.
$srcForCopy = $wb.Worksheets.Item(2)
$srcForCopy.Copy($wb.Worksheets.Item(1)) | Out-Null

$allPages = $wb.Worksheets.Item(1)
$allPages.Name = "All Pages"

# Clear out the copied sample data (rows 1-3, cols A-D) before writing the real content.
$allPages.Range("A1:D3").ClearContents() | Out-Null

# Row 1 - header
$allPages.Range("A1").Value = "Hostname"
$allPages.Range("B1").Value = "PageOptions Link"
$allPages.Range("C1").Value = "PageOptions Type"
$allPages.Range("D1").Value = "Comment"

# Row 2
$allPages.Range("A2").Value = "https://www-qa.cancer.gov"
$allPages.Range("B2").Value = "/about-cancer/understanding/what-is-cancer"
$allPages.Range("C2").Value = "all"
$allPages.Range("D2").Value = "all buttons, Content pages, article, factsheet, topic page, press release, institution, cancer research page, general page, pdq summaries, blog posts, blog series"
$allPages.Range("E2").Value = "FontResize"
$allPages.Range("F2").Value = "print"
$allPages.Range("G2").Value = "mail"
$allPages.Range("H2").Value = "facebook"
$allPages.Range("I2").Value = "twitter"
$allPages.Range("J2").Value = "Google+"
$allPages.Range("K2").Value = "Pintrest"

# Row 3
$allPages.Range("A3").Value = "https://www-qa.cancer.gov"
$allPages.Range("B3").Value = "/about-cancer"
$allPages.Range("C3").Value = "default"
$allPages.Range("D3").Value = "default buttons"
$allPages.Range("F3").Value = "print"
$allPages.Range("G3").Value = "mail"
$allPages.Range("H3").Value = "facebook"
$allPages.Range("I3").Value = "twitter"
$allPages.Range("J3").Value = "Google+"
$allPages.Range("K3").Value = "Pintrest"

# Row 4
$allPages.Range("A4").Value = "https://www-qa.cancer.gov"
$allPages.Range("B4").Value = "/research/resources"
$allPages.Range("C4").Value = "r4r"
$allPages.Range("D4").Value = "bottom button"
$allPages.Range("G4").Value = "mail"
$allPages.Range("I4").Value = "twitter"

# Row 5
$allPages.Range("A5").Value = "https://www-qa.cancer.gov"
$allPages.Range("B5").Value = "/espanol/cancer/naturaleza/que-es"
$allPages.Range("C5").Value = "all"
$allPages.Range("D5").Value = "all buttons, Content pages, article, factsheet, topic page, press release, institution, cancer research page, general page, pdq summaries, blog posts, blog series"
$allPages.Range("E5").Value = "FontResize"
$allPages.Range("F5").Value = "print"
$allPages.Range("G5").Value = "mail"
$allPages.Range("H5").Value = "facebook"
$allPages.Range("I5").Value = "twitter"
$allPages.Range("J5").Value = "Google+"
$allPages.Range("K5").Value = "Pintrest"

# Row 6
$allPages.Range("A6").Value = "https://www-qa.cancer.gov"
$allPages.Range("B6").Value = "/espanol/cancer"
$allPages.Range("C6").Value = "default"
$allPages.Range("D6").Value = "default buttons"
$allPages.Range("F6").Value = "print"
$allPages.Range("G6").Value = "mail"
$allPages.Range("H6").Value = "facebook"
$allPages.Range("I6").Value = "twitter"
$allPages.Range("J6").Value = "Google+"
$allPages.Range("K6").Value = "Pintrest"

# Row 7 (new)
$allPages.Range("A7").Value = "https://livehelp.cancer.gov"
$allPages.Range("B7").Value = "/"
$allPages.Range("C7").Value = "none"
$allPages.Range("D7").Value = "no buttons, CTS print pages, dictionary pop-ups, livehelp"

# Row 8 (new)
$allPages.Range("A8").Value = "https://livehelp-es.cancer.gov"
$allPages.Range("B8").Value = "/"
$allPages.Range("C8").Value = "none"
$allPages.Range("D8").Value = "no buttons, CTS print pages, dictionary pop-ups, livehelp"

# Selection / view for "All Pages": not the tabSelected sheet, cursor on A9
$allPages.Activate() | Out-Null
$allPages.Range("A9").Select() | Out-Null

# "Visible Page Options" sheet: keep tabSelected, just move the cursor to A7
$visible = $wb.Worksheets.Item("Visible Page Options")
$visible.Activate() | Out-Null
$visible.Range("A7").Select() | Out-Null

# "Not Visible Page Options" sheet: selection becomes A2:D3 (active cell A2)
$notVisible = $wb.Worksheets.Item("Not Visible Page Options")
$notVisible.Activate() | Out-Null
$notVisible.Range("A2:D3").Select() | Out-Null

# Final active sheet should be "Visible Page Options" (workbook activeTab = 1)
$visible.Activate() | Out-Null
